$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: add new trailing cells X3 (PriceChange) and Y3 (UpDown verdict "Up")
$ws.Range("X3").Value = 4.75
$ws.Range("Y3").Value = "Up"

# Row 4: new data row
$ws.Range("A4").Value = Get-Date -Year 2016 -Month 9 -Day 28 -Hour 21 -Minute 23 -Second 45
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = -4
$ws.Range("E4").Value = 10109
$ws.Range("F4").Value = 561
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 26
$ws.Range("J4").Value = 73
$ws.Range("K4").Value = 12293
$ws.Range("L4").Value = 118
$ws.Range("M4").Value = 52
$ws.Range("N4").Value = 4
$ws.Range("O4").Value = 11
$ws.Range("P4").Value = "Named"
$ws.Range("Q4").Value = 57.519894101767122
$ws.Range("R4").Value = 1.83
$ws.Range("S4").Value = 0.13639999999999999
$ws.Range("S4").NumberFormat = $ws.Range("S3").NumberFormat
$ws.Range("T4").Value = 0.016500000000000001
$ws.Range("T4").NumberFormat = $ws.Range("T3").NumberFormat
$ws.Range("U4").Value = 6.04
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 2
